# Update column F (numeric) values on each sheet per the target diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1506
$ws.Cells.Item(5, 6).Value = 7651
$ws.Cells.Item(6, 6).Value = 4835
$ws.Cells.Item(7, 6).Value = 7120
$ws.Cells.Item(10, 6).Value = 1511
$ws.Cells.Item(11, 6).Value = 873
$ws.Cells.Item(16, 6).Value = 179
$ws.Cells.Item(20, 6).Value = 237
$ws.Cells.Item(27, 6).Value = 50
$ws.Cells.Item(28, 6).Value = 155
$ws.Cells.Item(33, 6).Value = 6
$ws.Cells.Item(35, 6).Value = 117
$ws.Cells.Item(36, 6).Value = 38
$ws.Cells.Item(38, 6).Value = 428
$ws.Cells.Item(41, 6).Value = 94
$ws.Cells.Item(42, 6).Value = 405
$ws.Cells.Item(43, 6).Value = 1204
$ws.Cells.Item(44, 6).Value = 589
$ws.Cells.Item(46, 6).Value = 26

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(21, 6).Value = 9
$ws.Cells.Item(24, 6).Value = 142
$ws.Cells.Item(27, 6).Value = 641
$ws.Cells.Item(31, 6).Value = 116
$ws.Cells.Item(32, 6).Value = 868
$ws.Cells.Item(45, 6).Value = 79

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(5, 6).Value = 857
$ws.Cells.Item(8, 6).Value = 69
$ws.Cells.Item(9, 6).Value = 73
$ws.Cells.Item(10, 6).Value = 1670
$ws.Cells.Item(11, 6).Value = 2573

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 1506
$ws.Cells.Item(4, 6).Value = 857
$ws.Cells.Item(9, 6).Value = 7651
$ws.Cells.Item(10, 6).Value = 4835
$ws.Cells.Item(11, 6).Value = 7120
$ws.Cells.Item(13, 6).Value = 1511
$ws.Cells.Item(15, 6).Value = 873
$ws.Cells.Item(18, 6).Value = 1670
$ws.Cells.Item(19, 6).Value = 2573
$ws.Cells.Item(23, 6).Value = 179
$ws.Cells.Item(25, 6).Value = 237
$ws.Cells.Item(27, 6).Value = 641
$ws.Cells.Item(30, 6).Value = 155
$ws.Cells.Item(33, 6).Value = 868
$ws.Cells.Item(35, 6).Value = 117
$ws.Cells.Item(41, 6).Value = 94
$ws.Cells.Item(43, 6).Value = 405
$ws.Cells.Item(44, 6).Value = 589
$ws.Cells.Item(49, 6).Value = 79
$ws.Cells.Item(50, 6).Value = 26

